$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "... wird auf den L|e|benszyklus ..." -> merge the split runs into one
#    continuous run of text (no formatting change, just text reflow).
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    " wird auf den Lebenszyklus von HTML Elementen eingegangen. Vor allem werden wir sehen, dass ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " wird auf den Lebenszyklus von HTML Elementen eingegangen. Vor allem werden wir sehen, dass ",
    2)

# ---------------------------------------------------------------------------
# 2) "bei gle|(bookmark _GoBack)|ichem Endergebnis." -> merge text, and the
#    _GoBack bookmark that used to sit here disappears (it moves further
#    down in the document, see step below).
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "bei gleichem Endergebnis.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "bei gleichem Endergebnis.",
    2)

# ---------------------------------------------------------------------------
# 3) " ke|n|nen lernen." -> merge into " kennen lernen."
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    " kennen lernen.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " kennen lernen.",
    2)

# ---------------------------------------------------------------------------
# 4) Replace the bulleted outline under "Schiffe versenken"
#    ("Spielidee vorstellen" ... "(Übung) Erneutes Verstecken ohne Browser
#    Refresh", 9 numbered list paragraphs) with 6 plain narrative paragraphs.
# ---------------------------------------------------------------------------

function Find-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# Bold a single (first, from the start of $para's range) occurrence of $word
# inside paragraph $para.
function Set-BoldWord($para, $needle) {
    $rng = $para.Range.Duplicate
    $ok = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $rng.Font.Bold = 1
    }
}

# Italicize a single (first) occurrence of $needle inside paragraph $para.
function Set-ItalicText($para, $needle) {
    $rng = $para.Range.Duplicate
    $ok = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $rng.Font.Italic = 1
    }
}

$firstOld = Find-ParagraphByText $d "Spielidee vorstellen"
$lastOld = Find-ParagraphByText $d "(Übung) Erneutes Verstecken ohne Browser Refresh"

$delRange = $d.Range($firstOld.Range.Start, $lastOld.Range.End)
$delRange.Delete()

$heading = Find-ParagraphByText $d "Schiffe versenken"

$newTexts = @(
    "Der letzte (optionale) Abschnitt stellt eine andere Alternative zum Aufbau eines Spielfelds vor und verwendet ansonsten die Grundlagen, die in den vorherigen Abschnitten bereits erarbeitet wurden. Als Grundlage dient das wohlbekannte Spiel Schiffe versenken – natürlich in der Umsetzung stark vereinfacht.",
    "In der ersten Übung wird wieder ein Spielfeld aufgebaut. Da dessen Größe aber mit 10x10 bekannt ist, wird dieses hier statisch in HTML auf Basis von DIV Tabellen aufgebaut. Wir lernen dabei die HTML Entität &nbsp; sowie die CSS Eigenschaften text-align und vertical-align kennen. Nachdem wir uns Gedanken über die Anzeige verschiedener Zustände der Spielzellen gemacht haben (bereits beschossen oder nicht, Treffer oder nicht) nutzen wird CSS Klassen zum Aufbau einer Legende. Im JavaScript Programmcode prüfen wir die Korrektheit des Spielfelds – statt es dynamisch aufzubauen. Wir lernen dabei die Auswahl von HTML Elementen mit :not kennen.",
    "Diesmal ist bereits die zweite Übung eine echte Herausforderung und es gibt eine Vielzahl von Möglichkeiten der Umsetzung. Ziel der Übung ist es, dass der erste Spieler seine Schiffe auf dem Spielfeld versteckt – wie gewohnt ein 5er, ein 4er, zwei 3er und ein 2er. Die Beispiellösung der Übung verwendet das Ziehen und Fallenlassen mit der Maus um diese Aufgabe zu meistern. Wichtig ist es dabei zu beachten, dass jedes Schiff in zwei Orientierungen (längs und quer) positioniert werden kann. Je nach Orientierung natürlich nicht an jeder Stelle auf dem Spielfeld. Im Ansatz mit Ziehen und Fallenlassen muss während des Ziehens bereits geprüft werden, ob ein Schiff in einer bestimmten Orientierung fallengelassen werden kann. Erst das Fallenlassen markiert dann die entsprechenden Spielzellen. Zusätzlich soll es dem ersten Spieler möglich sein, eine einmal getroffene Entscheidung zu revidieren und ein Schiff an eine andere Stelle zu versetzen – möglicherweise sogar mit einer anderen Orientierung. Solange das Spielfeld vom ersten Spieler aufgebaut wird muss dieser natürlich seine Auswahl sehen können. Ist er mit seinen Verstecken zufrieden, schaltet er auf den Suchmodus für den zweiten Spieler um.",
    "Die Suche durch den zweiten Spieler wird dann die dritte Übung sein. Beim Anklicken einer Spielzelle muss auf einen möglichen Treffer geprüft werden. Wenn es einen Treffer gibt, soll weiterhin untersucht werden, ob ein Schiff als Ganzes versenkt wurde. Mit dem Versenken des letzten Schiffs ist das Spiel beendet und der zweite Spieler soll nun die Anzahl der Versuche angezeigt bekommen.  Wir werden hier mit JavaScript Objekten arbeiten, die an die einzelnen HTML Elemente der Spielzellen gebunden werden. Einen Treffer zu erkennen wird dann sehr einfach – versenkte Schiffe sind schon etwas kniffeliger.",
    "Damit das Spiel nicht zu schwer und frustrierend wird, bauen wir in der vierten Übung die Option ein, sich bei jedem Versuch anzeigen zu lassen, wie viele Treffer sich in den (3 bis 8) umliegenden Zellen ergeben. Als praktisch erweisen sich dabei die Funktionen min und max aus der Math Bibliothek.",
    "Die letzte Übung erlaubt einen Neustart des Spiels, ohne einen Browser Refresh durchführen zu müssen. Hier können wir lernen, wie man die im Spielverlauf durchgeführten Veränderungen an den HTML Elementen ordnungsgemäß rückgängig macht – in den meisten Fällen geht es um die className Eigenschaft sowie eigene Erweiterungen in Form von angehängten JavaScript Objekten."
)

$prev = $heading
$newParas = @()
foreach ($txt in $newTexts) {
    $prev.Range.InsertParagraphAfter()
    $cur = $prev.Next()
    $cur.Style = "Normal"
    $cur.Range.Text = $txt
    $newParas += $cur
    $prev = $cur
}

# Bold the (first) "Übung" in each of the 5 paragraphs that mention it.
Set-BoldWord $newParas[1] "Übung"
Set-BoldWord $newParas[2] "Übung"
Set-BoldWord $newParas[3] "Übung"
Set-BoldWord $newParas[4] "Übung"
Set-BoldWord $newParas[5] "Übung"

# Italics.
Set-ItalicText $newParas[0] "Schiffe versenken"
Set-ItalicText $newParas[1] "&nbsp;"
Set-ItalicText $newParas[1] "text-align"
Set-ItalicText $newParas[1] "vertical-align"
Set-ItalicText $newParas[1] ":not"
Set-ItalicText $newParas[4] "min"
Set-ItalicText $newParas[4] "max"
Set-ItalicText $newParas[4] "Math"
Set-ItalicText $newParas[5] "className"

# The "_GoBack" bookmark that disappeared in step 2) re-appears around
# "Math" ("... aus der |Math |Bibliothek.").
$mathRng = $newParas[4].Range.Duplicate
$null = $mathRng.Find.Execute("Math", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmStart = $mathRng.Start       # right before "Math"
$bmEnd = $mathRng.End + 1       # include the trailing space before "Bibliothek."
$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
